$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 21:22"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1079894
$ws.Range("C4").Value = 15700
$ws.Range("E4").Value = 867536
$ws.Range("G4").Value = 1017
$ws.Range("H4").Value = 62672

# Francia (row 8)
$ws.Range("B8").Value = 167178
$ws.Range("C8").Value = 758
$ws.Range("E8").Value = 93326

# Alemania (row 9)
$ws.Range("B9").Value = 162375
$ws.Range("C9").Value = 836
$ws.Range("E9").Value = 32312
$ws.Range("G9").Value = 96
$ws.Range("H9").Value = 6563

# Peru overtakes India: row 18 now Peru (new higher numbers), row 19 now India (old India numbers)
$ws.Range("A18").Value = "Peru"
$ws.Range("B18").Value = 36976
$ws.Range("C18").Value = 3045
$ws.Range("D18").Value = 10405
$ws.Range("E18").Value = 25520
$ws.Range("F18").Value = 651
$ws.Range("G18").Value = 108
$ws.Range("H18").Value = 1051

$ws.Range("A19").Value = "India"
$ws.Range("B19").Value = 34780
$ws.Range("C19").Value = 1718
$ws.Range("D19").Value = 9068
$ws.Range("E19").Value = 24561
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 72
$ws.Range("H19").Value = 1151

# Pakistan overtakes Singapur: row 27 now Pakistan (new higher numbers), row 28 now Singapur (old Singapur numbers)
$ws.Range("A27").Value = "Pakistan"
$ws.Range("B27").Value = 16473
$ws.Range("C27").Value = 948
$ws.Range("D27").Value = 4105
$ws.Range("E27").Value = 12007
$ws.Range("F27").Value = 111
$ws.Range("G27").Value = 18
$ws.Range("H27").Value = 361

$ws.Range("A28").Value = "Singapur"
$ws.Range("B28").Value = 16169
$ws.Range("C28").Value = 528
$ws.Range("D28").Value = 1244
$ws.Range("E28").Value = 14910
$ws.Range("F28").Value = 22
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 15

# Argentina (row 56)
$ws.Range("D56").Value = 1256
$ws.Range("E56").Value = 2813
$ws.Range("F56").Value = 157
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 216

# Barein (row 62)
$ws.Range("B62").Value = 3040
$ws.Range("C62").Value = 119
$ws.Range("D62").Value = 1500
$ws.Range("E62").Value = 1532
